$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 32
$ws1.Range("F4").Value = 197
$ws1.Range("F5").Value = 1109
$ws1.Range("F6").Value = 8286
$ws1.Range("F7").Value = 8286
$ws1.Range("F8").Value = 141
$ws1.Range("F10").Value = 6930
$ws1.Range("F12").Value = 5060
$ws1.Range("F13").Value = 5531
$ws1.Range("F14").Value = 1079
$ws1.Range("F15").Value = 341
$ws1.Range("F16").Value = 349
$ws1.Range("F25").Value = 9292
$ws1.Range("F26").Value = 74
$ws1.Range("F27").Value = 1704
$ws1.Range("F28").Value = 982
$ws1.Range("F31").Value = 1895
$ws1.Range("F37").Value = 1908
$ws1.Range("F38").Value = 246
$ws1.Range("F39").Value = 1215
$ws1.Range("F41").Value = 4851
$ws1.Range("F42").Value = 379
$ws1.Range("F43").Value = 1166
$ws1.Range("F44").Value = 81
$ws1.Range("F47").Value = 1082
$ws1.Range("F48").Value = 1046
$ws1.Range("F49").Value = 928
$ws1.Range("F50").Value = 1275
$ws1.Range("F51").Value = 46

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 42
$ws2.Range("F4").Value = 1
$ws2.Range("F5").Value = 2
$ws2.Range("F9").Value = 6
$ws2.Range("F16").Value = 89

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 32
$ws4.Range("F4").Value = 197
$ws4.Range("F5").Value = 42
$ws4.Range("F6").Value = 1109
$ws4.Range("F7").Value = 8286
$ws4.Range("F8").Value = 141
$ws4.Range("F10").Value = 6930
$ws4.Range("F14").Value = 5060
$ws4.Range("F15").Value = 5531
$ws4.Range("F16").Value = 1079
$ws4.Range("F17").Value = 341
$ws4.Range("F18").Value = 349
$ws4.Range("F25").Value = 9292
$ws4.Range("F26").Value = 74
$ws4.Range("F27").Value = 1704
$ws4.Range("F28").Value = 982
$ws4.Range("F31").Value = 1895
$ws4.Range("F37").Value = 1908
$ws4.Range("F38").Value = 246
$ws4.Range("F39").Value = 1215
$ws4.Range("F41").Value = 4851
$ws4.Range("F42").Value = 379
$ws4.Range("F43").Value = 1166
$ws4.Range("F44").Value = 81
$ws4.Range("F47").Value = 1082
$ws4.Range("F48").Value = 1046
$ws4.Range("F49").Value = 928
$ws4.Range("F50").Value = 1275
$ws4.Range("F51").Value = 46
